$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEGFA165_VEGFR1")

# The sheet currently ends at row 8 (last row has the thick bottom border).
# We need to add a new last row (row 9) for "Breier et al., 1995" and turn the
# old last row (row 8) into a normal interior row.

# 1) Push row 8's current ("last row") formatting down onto the new row 9,
#    so row 9 becomes the new bottom-bordered row.
$ws.Range("A8:D8").Copy() | Out-Null
$ws.Range("A9:D9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 2) Restore row 8 to a normal interior-row look by copying row 7's format
#    (no thick bottom border) onto it.
$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A8:D8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Rows(8).AutoFit()

$ws.Application.CutCopyMode = $false

# 3) Fill in the new reference row.
$ws.Range("A9").Value = "Breier et al., 1995"
$ws.Range("B9").Value = "Radioligand"
$ws.Range("C9").Value = 114
$ws.Range("D9").Value = ""

# 4) Excel's style table gets compacted/deduplicated whenever the workbook is
#    resaved after this kind of edit, which also normalizes a couple of
#    vestigial duplicate "applyFill" border styles elsewhere in the workbook
#    (these carry no actual fill) back onto their canonical equivalents.
#    Reproduce that on the untouched VEGFA165_VEGFR2 sheet's closing row.
$ws3 = $wb.Worksheets.Item("VEGFA165_NRP1")
$ws2 = $wb.Worksheets.Item("VEGFA165_VEGFR2")
$ws3.Range("A7:D7").Copy() | Out-Null
$ws2.Range("A11:D11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# 5) Make this sheet the active/selected tab, matching the author's final view.
$ws.Activate()
$ws.Range("C17").Select() | Out-Null
